# Updates coin price / 1h-volume figures (and two swapped-row fixes) in the
# crypto ranking sheet, matching the refreshed data from the source feed.
#
# Columns D (Price) and E (Volume 1h) are stored as literal text in this sheet
# (e.g. "60.929.54" uses dots as thousand separators, and "1.00" / "0.999" must
# stay text, not become the number 1). Assigning a numeric-looking string via
# .Value normally gets auto-coerced to a real number by Excel, so for any new
# value that Excel would parse as a number we briefly force the cell to Text
# format, assign the value, then restore the "Normal" style so the cell
# formatting stays exactly as it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.987.40"
$ws.Range("E2").Value = "  -3.45%  "

# Row 3
$ws.Range("D3").Value = "3.357.61"
$ws.Range("E3").Value = "  -2.84%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.70%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.54%  "

# Row 7
$ws.Range("E7").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.98"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.84%  "

# Row 10
$ws.Range("E10").Value = "  -1.26%  "

# Row 11
$ws.Range("E11").Value = "  +1.86%  "

# Row 12
$ws.Range("D12").Value = "3.934.84"
$ws.Range("E12").Value = "  -2.79%  "

# Row 13
$ws.Range("E13").Value = "  +1.01%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.95%  "

# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000170"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.11%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.353.67"
$ws.Range("E16").Value = "  -3.28%  "

# Row 17
$ws.Range("D17").Value = "61.074.35"
$ws.Range("E17").Value = "  -3.27%  "

# Row 18
$ws.Range("E18").Value = "  -1.10%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.55%  "

# Row 20
$ws.Range("E20").Value = "  -3.16%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.34%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.563"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.35%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "75.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.08%  "

# Row 24
$ws.Range("E24").Value = "  +0.03%  "

# Row 25
$ws.Range("D25").Value = "3.525.85"
$ws.Range("E25").Value = "  -1.76%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000109"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.15%  "

# Row 27
$ws.Range("E27").Value = "  -3.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.25%  "

# Row 29
$ws.Range("E29").Value = "  -3.66%  "

# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.78%  "

# Row 31
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.02%  "

# Row 32
$ws.Range("E32").Value = "  -4.71%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.86%  "

# Row 34
$ws.Range("E34").Value = "  -3.99%  "

# Row 35
$ws.Range("E35").Value = "  +0.32%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "169.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.20%  "

# Row 37
$ws.Range("E37").Value = "  -4.40%  "

# Row 38
$ws.Range("E38").Value = "  -4.07%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "29.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.56%  "

# Row 40
$ws.Range("D40").Value = "3.394.17"
$ws.Range("E40").Value = "  -2.76%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0756"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.59%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.29%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.762"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.62%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.54%  "

# Row 45
$ws.Range("E45").Value = "  -3.90%  "

# Row 46
$ws.Range("E46").Value = "  -6.35%  "

# Row 47
$ws.Range("D47").Value = "2.511.72"
$ws.Range("E47").Value = "  -2.70%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.92%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.24%  "

# Row 50
$ws.Range("E50").Value = "  +0.05%  "

# Row 51
$ws.Range("E51").Value = "  -2.33%  "
